# Auto-generated: updates FFXIV Leve profit market-data columns (H:N)
# on each class sheet to match the scheduled market-data refresh.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 96 (Leve Item ID 19894)
$ws.Range("H96").Value = 135.3
$ws.Range("I96").Value = 139.22223
$ws.Range("J96").Value = 100
$ws.Range("K96").Value = 417.66669
$ws.Range("L96").Value = 300
$ws.Range("M96").Value = 955.33331
$ws.Range("N96").Value = -3046
# Row 104 (Leve Item ID 24263)
$ws.Range("H104").Value = 154.33333
$ws.Range("I104").Value = 154.33333
$ws.Range("K104").Value = 462.99999
$ws.Range("M104").Value = 1284.00001
# Row 111 (Leve Item ID 27768)
$ws.Range("H111").Value = 5994.6
$ws.Range("I111").Value = 4831.364
$ws.Range("J111").Value = 9193.5
$ws.Range("K111").Value = 14494.092
$ws.Range("L111").Value = 27580.5
$ws.Range("M111").Value = -11427.092
$ws.Range("N111").Value = -33714.5
# Row 137 (Leve Item ID 44013)
$ws.Range("H137").Value = 1595.7576
$ws.Range("I137").Value = 1125.5
$ws.Range("K137").Value = 3376.5
$ws.Range("M137").Value = -826.5

$ws = $wb.Worksheets.Item("ARM")
# Row 61 (Leve Item ID 43999)
$ws.Range("H61").Value = 4379.875
$ws.Range("I61").Value = 1933.1666
$ws.Range("J61").Value = 5847.9
$ws.Range("K61").Value = 1933.1666
$ws.Range("L61").Value = 5847.9
$ws.Range("M61").Value = -1721.1666
$ws.Range("N61").Value = -6271.9
# Row 122 (Leve Item ID 36168)
$ws.Range("H122").Value = 3875.0557
$ws.Range("I122").Value = 2491.4348
$ws.Range("K122").Value = 7474.3044
$ws.Range("M122").Value = -5024.3044
# Row 132 (Leve Item ID 43997)
$ws.Range("H132").Value = 3566.5833
$ws.Range("I132").Value = 2899.8572
$ws.Range("J132").Value = 4500
$ws.Range("K132").Value = 8699.571599999999
$ws.Range("L132").Value = 13500
$ws.Range("M132").Value = -6169.571599999999
$ws.Range("N132").Value = -18560
# Row 136 (Leve Item ID 43999)
$ws.Range("H136").Value = 4379.875
$ws.Range("I136").Value = 1933.1666
$ws.Range("J136").Value = 5847.9
$ws.Range("K136").Value = 5799.4998
$ws.Range("L136").Value = 17543.7
$ws.Range("M136").Value = -3249.4998
$ws.Range("N136").Value = -22643.7

$ws = $wb.Worksheets.Item("BSM")
# Row 20 (Leve Item ID 14149)
$ws.Range("H20").Value = 4433
$ws.Range("I20").Value = 2574.5833
$ws.Range("K20").Value = 2574.5833
$ws.Range("M20").Value = -2327.5833
# Row 99 (Leve Item ID 19943)
$ws.Range("H99").Value = 23548.947
$ws.Range("I99").Value = 27250.562
$ws.Range("J99").Value = 3807
$ws.Range("K99").Value = 27250.562
$ws.Range("L99").Value = 3807
$ws.Range("M99").Value = -25752.562
$ws.Range("N99").Value = -6803
# Row 105 (Leve Item ID 19947)
$ws.Range("H105").Value = 2009.8
$ws.Range("I105").Value = 2109.7778
$ws.Range("K105").Value = 2109.7778
$ws.Range("M105").Value = -362.7777999999998
# Row 111 (Leve Item ID 25789)
$ws.Range("H111").Value = 20000
$ws.Range("J111").Value = 20000
$ws.Range("L111").Value = 20000
$ws.Range("N111").Value = -28180
# Row 119 (Leve Item ID 26281)
$ws.Range("H119").Value = 0
$ws.Range("J119").Value = 0
$ws.Range("L119").Value = 0
$ws.Range("N119").ClearContents()
# Row 134 (Leve Item ID 43998)
$ws.Range("H134").Value = 2716.4348
$ws.Range("I134").Value = 1942.4872
$ws.Range("J134").Value = 7028.4287
$ws.Range("K134").Value = 5827.461600000001
$ws.Range("L134").Value = 21085.2861
$ws.Range("M134").Value = -3292.461600000001
$ws.Range("N134").Value = -26155.2861

$ws = $wb.Worksheets.Item("CRP")
# Row 62 (Leve Item ID 12580)
$ws.Range("H62").Value = 6572.273
$ws.Range("J62").Value = 4150
$ws.Range("L62").Value = 4150
$ws.Range("N62").Value = -5398
# Row 65 (Leve Item ID 12580)
$ws.Range("H65").Value = 6572.273
$ws.Range("J65").Value = 4150
$ws.Range("L65").Value = 20750
$ws.Range("N65").Value = -26990
# Row 93 (Leve Item ID 19516)
$ws.Range("H93").Value = 29374.75
$ws.Range("I93").Value = 5833
$ws.Range("J93").Value = 100000
$ws.Range("K93").Value = 5833
$ws.Range("L93").Value = 100000
$ws.Range("M93").Value = -3961
$ws.Range("N93").Value = -103744
# Row 103 (Leve Item ID 19558)
$ws.Range("H103").Value = 10000
$ws.Range("I103").Value = 10000
$ws.Range("K103").Value = 10000
$ws.Range("M103").Value = -8828
# Row 120 (Leve Item ID 27230)
$ws.Range("H120").Value = 46497.5
$ws.Range("I120").Value = 0
$ws.Range("J120").Value = 46497.5
$ws.Range("K120").Value = 0
$ws.Range("L120").Value = 46497.5
$ws.Range("M120").ClearContents()
$ws.Range("N120").Value = -53755.5
# Row 134 (Leve Item ID 44020)
$ws.Range("H134").Value = 2904.3833
$ws.Range("I134").Value = 1658.6904
$ws.Range("K134").Value = 4976.0712
$ws.Range("M134").Value = -2441.0712

$ws = $wb.Worksheets.Item("CUL")
# Row 4 (Leve Item ID 4650)
$ws.Range("H4").Value = 33371028
$ws.Range("I4").Value = 37078810
$ws.Range("J4").Value = 998.6667
$ws.Range("K4").Value = 111236430
$ws.Range("L4").Value = 2996.0001
$ws.Range("M4").Value = -111236318
$ws.Range("N4").Value = -3220.0001
# Row 5 (Leve Item ID 43974)
$ws.Range("H5").Value = 7511.2583
$ws.Range("I5").Value = 4042
$ws.Range("J5").Value = 11723.929
$ws.Range("K5").Value = 12126
$ws.Range("L5").Value = 35171.787
$ws.Range("M5").Value = -12014
$ws.Range("N5").Value = -35395.787
# Row 38 (Leve Item ID 4860)
$ws.Range("H38").Value = 2097.2
$ws.Range("I38").Value = 89
$ws.Range("J38").Value = 5109.5
$ws.Range("K38").Value = 267
$ws.Range("L38").Value = 15328.5
$ws.Range("M38").Value = 80
$ws.Range("N38").Value = -16022.5
# Row 40 (Leve Item ID 4827)
$ws.Range("H40").Value = 56.77778
$ws.Range("J40").Value = 50.25
$ws.Range("L40").Value = 201
$ws.Range("N40").Value = -339
# Row 86 (Leve Item ID 12892)
$ws.Range("H86").Value = 793.0833
$ws.Range("J86").Value = 809.2
$ws.Range("L86").Value = 2427.6
$ws.Range("N86").Value = -4799.6
# Row 89 (Leve Item ID 12892)
$ws.Range("H89").Value = 793.0833
$ws.Range("J89").Value = 809.2
$ws.Range("L89").Value = 7282.8
$ws.Range("N89").Value = -19138.8
# Row 113 (Leve Item ID 27843)
$ws.Range("H113").Value = 6458.1055
$ws.Range("I113").Value = 15314
$ws.Range("K113").Value = 45942
$ws.Range("M113").Value = -43772
# Row 115 (Leve Item ID 27861)
$ws.Range("H115").Value = 3113.7646
$ws.Range("I115").Value = 1379.25
$ws.Range("J115").Value = 4655.5557
$ws.Range("K115").Value = 4137.75
$ws.Range("L115").Value = 13966.6671
$ws.Range("M115").Value = -2962.75
$ws.Range("N115").Value = -16316.6671
# Row 122 (Leve Item ID 36078)
$ws.Range("H122").Value = 1186.4546
$ws.Range("J122").Value = 1219.3529
$ws.Range("L122").Value = 10974.1761
$ws.Range("N122").Value = -15874.1761
# Row 123 (Leve Item ID 36037)
$ws.Range("H123").Value = 4105.8
$ws.Range("I123").Value = 2030
$ws.Range("K123").Value = 6090
$ws.Range("M123").Value = -3640
# Row 131 (Leve Item ID 36060)
$ws.Range("H131").Value = 4710.55
$ws.Range("I131").Value = 2813.5454
$ws.Range("J131").Value = 7029.1113
$ws.Range("K131").Value = 8440.636200000001
$ws.Range("L131").Value = 21087.3339
$ws.Range("M131").Value = -3400.636200000001
$ws.Range("N131").Value = -31167.3339
# Row 135 (Leve Item ID 43974)
$ws.Range("H135").Value = 7511.2583
$ws.Range("I135").Value = 4042
$ws.Range("J135").Value = 11723.929
$ws.Range("K135").Value = 36378
$ws.Range("L135").Value = 105515.361
$ws.Range("M135").Value = -33843
$ws.Range("N135").Value = -110585.361

$ws = $wb.Worksheets.Item("GSM")
# Row 97 (Leve Item ID 19940)
$ws.Range("H97").Value = 1010
$ws.Range("I97").Value = 1010
$ws.Range("J97").Value = 0
$ws.Range("K97").Value = 1010
$ws.Range("L97").Value = 0
$ws.Range("M97").Value = -514
$ws.Range("N97").ClearContents()
# Row 132 (Leve Item ID 44008)
$ws.Range("H132").Value = 3139.7693
$ws.Range("I132").Value = 1657.1177
$ws.Range("K132").Value = 4971.3531
$ws.Range("M132").Value = -2441.3531

$ws = $wb.Worksheets.Item("LTW")
# Row 7 (Leve Item ID 36249)
$ws.Range("H7").Value = 3772.2727
$ws.Range("I7").Value = 2225.7693
$ws.Range("J7").Value = 6006.1113
$ws.Range("K7").Value = 2225.7693
$ws.Range("L7").Value = 6006.1113
$ws.Range("M7").Value = -2113.7693
$ws.Range("N7").Value = -6230.1113
# Row 22 (Leve Item ID 5277)
$ws.Range("H22").Value = 764.0833
$ws.Range("J22").Value = 999
$ws.Range("L22").Value = 999
$ws.Range("N22").Value = -1589
# Row 27 (Leve Item ID 5277)
$ws.Range("H27").Value = 764.0833
$ws.Range("J27").Value = 999
$ws.Range("L27").Value = 999
$ws.Range("N27").Value = -1213
# Row 40 (Leve Item ID 36248)
$ws.Range("H40").Value = 11197.556
$ws.Range("I40").Value = 14893.875
$ws.Range("J40").Value = 8240.5
$ws.Range("K40").Value = 14893.875
$ws.Range("L40").Value = 8240.5
$ws.Range("M40").Value = -14757.875
$ws.Range("N40").Value = -8512.5
# Row 46 (Leve Item ID 5282)
$ws.Range("H46").Value = 2375
$ws.Range("I46").Value = 2100
$ws.Range("J46").Value = 2650
$ws.Range("K46").Value = 2100
$ws.Range("L46").Value = 2650
$ws.Range("M46").Value = -1912
$ws.Range("N46").Value = -3026
# Row 126 (Leve Item ID 36249)
$ws.Range("H126").Value = 3772.2727
$ws.Range("I126").Value = 2225.7693
$ws.Range("J126").Value = 6006.1113
$ws.Range("K126").Value = 6677.3079
$ws.Range("L126").Value = 18018.3339
$ws.Range("M126").Value = -4207.3079
$ws.Range("N126").Value = -22958.3339
# Row 127 (Leve Item ID 34401)
$ws.Range("H127").Value = 0
$ws.Range("J127").Value = 0
$ws.Range("L127").Value = 0
$ws.Range("N127").ClearContents()
# Row 132 (Leve Item ID 44058)
$ws.Range("H132").Value = 3128.195
$ws.Range("I132").Value = 2543.4
$ws.Range("J132").Value = 4723.091
$ws.Range("K132").Value = 7630.200000000001
$ws.Range("L132").Value = 14169.273
$ws.Range("M132").Value = -5100.200000000001
$ws.Range("N132").Value = -19229.273
# Row 136 (Leve Item ID 44060)
$ws.Range("H136").Value = 4751.4
$ws.Range("I136").Value = 3015.4211
$ws.Range("K136").Value = 9046.263300000001
$ws.Range("M136").Value = -6496.263300000001

$ws = $wb.Worksheets.Item("WVR")
# Row 25 (Leve Item ID 3064)
$ws.Range("H25").Value = 0
$ws.Range("J25").Value = 0
$ws.Range("L25").Value = 0
$ws.Range("N25").ClearContents()
# Row 81 (Leve Item ID 12596)
$ws.Range("H81").Value = 7026
$ws.Range("I81").Value = 11779.4
$ws.Range("J81").Value = 1744.4445
$ws.Range("K81").Value = 23558.8
$ws.Range("L81").Value = 3488.889
$ws.Range("M81").Value = -22497.8
$ws.Range("N81").Value = -5610.889
# Row 84 (Leve Item ID 12596)
$ws.Range("H84").Value = 7026
$ws.Range("I84").Value = 11779.4
$ws.Range("J84").Value = 1744.4445
$ws.Range("K84").Value = 117794
$ws.Range("L84").Value = 17444.445
$ws.Range("M84").Value = -112490
$ws.Range("N84").Value = -28052.445
# Row 100 (Leve Item ID 19981)
$ws.Range("H100").Value = 658.3043
$ws.Range("I100").Value = 580.7895
$ws.Range("J100").Value = 1026.5
$ws.Range("K100").Value = 1161.579
$ws.Range("L100").Value = 2053
$ws.Range("M100").Value = -620.579
$ws.Range("N100").Value = -3135
# Row 132 (Leve Item ID 44029)
$ws.Range("H132").Value = 1680.8246
$ws.Range("I132").Value = 1059.4681
$ws.Range("K132").Value = 3178.4043
$ws.Range("M132").Value = -648.4043000000001
# Row 141 (Leve Item ID 42505)
$ws.Range("H141").Value = 99999.234
$ws.Range("J141").Value = 99999.234
$ws.Range("L141").Value = 99999.234
$ws.Range("N141").Value = -110359.234

